$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text format so numeric-looking strings
# (e.g. "311.08") are not auto-converted to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "46.204.96"
$ws.Range("E2").Value = "  -1.80%  "
$ws.Range("D3").Value = "2.656.48"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "311.08"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "99.05"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("D7").Value = "0.598"
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "38.99"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "54.41"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "8.11"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("D14").Value = "3.033.05"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "2.643.64"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "0.926"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "14.97"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "46.195.06"
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "6.81"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "12.85"
$ws.Range("E22").Value = "  -4.06%  "
$ws.Range("D23").Value = "75.00"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("D24").Value = "284.32"
$ws.Range("E24").Value = "  +8.37%  "
$ws.Range("D25").Value = "3.06"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "30.28"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "10.60"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").Value = "38.95"
$ws.Range("E30").Value = "  -6.05%  "
$ws.Range("D31").Value = "2.25"
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("D32").Value = "6.28"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").Value = "3.76"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("D34").Value = "2.37"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").Value = "156.87"
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("D36").Value = "0.0845"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("E38").Value = "  +3.65%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "23.10"
$ws.Range("E40").Value = "  +6.77%  "
$ws.Range("D41").Value = "15.91"
$ws.Range("E41").Value = "  -6.48%  "
$ws.Range("D42").Value = "0.0330"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "3.60"
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("D44").Value = "4.06"
$ws.Range("E44").Value = "  -6.92%  "
$ws.Range("D45").Value = "2.154.61"
$ws.Range("E45").Value = "  +3.76%  "
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "94.55"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "111.07"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("D49").Value = "9.17"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.202"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.883.29"
$ws.Range("E51").Value = "  -0.74%  "

# Restore the original (default) cell style now that the text values are set,
# so no stray number-format style is left behind on the cells.
$dRange.Style = "Normal"

